$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("B2").Value = "1L-1AP-4AYB"
$ws.Range("C2").Value = "3-3-1"
$ws.Range("D2").Value = 29084
$ws.Range("E2").Value = 45815.82576238426

# --- Data for new rows 3..12 ---
$data = @(
    @(2,  "1P-1AYB-2AYB",              "3-1-2",     7045,   45814.82576238426),
    @(3,  "1AYB-1L-4AYB-1AP",          "1-1-1-3",   46501,  45812.82576238426),
    @(4,  "1AP-1M-5AYB-1P-1L",         "3-2-3-3-3", 33424,  45810.82576238426),
    @(5,  "1L-1AP",                    "2-1",       15429,  45805.82576238426),
    @(6,  "1AYB-1P",                   "3-1",       23108,  45795.82576238426),
    @(7,  "5AYB-1AP",                  "2-3",       25615,  45770.82576238426),
    @(8,  "3AYB-1M",                   "2-3",       43601,  45755.82576238426),
    @(9,  "1AP-3AYB-1AYB-5AYB-1M",     "1-3-1-3-2", 32203,  45809.82576238426),
    @(10, "1L-5AYB-1M",                "2-2-3",     31254,  45793.82576238426),
    @(11, "1AYB-1L-2M",                "4-4-6",     246000, 45815.82793382008)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r++
}
